# Automation script for testcase2
# Adds a new "InvalidLogin" worksheet after the existing "ValidLogin" sheet,
# containing a username/password header row and an invalid credential pair.

$wb = $excel.ActiveWorkbook

# Insert the new sheet right after the last existing sheet (ValidLogin),
# which makes it the new active/selected sheet - matching the target
# workbook's sheet order and activeTab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "InvalidLogin"

# Header row
$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "password"

# Invalid credentials row
$ws.Range("A2").Value = "abcd"
$ws.Range("B2").Value = "xyz"

# Leave selection on B2, as in the target workbook
$ws.Range("B2").Select()
